$d = $word.ActiveDocument

# Locate the target paragraph: "Contributed 50+ hours to community tutoring..."
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Contributed 50+ hours to community tutoring*") {
        $target = $p
    }
}
if ($target -eq $null) {
    throw "Could not find target paragraph"
}
$startIdx = $target.Index

# Create 6 placeholder paragraphs right after the target paragraph in one shot
# (using carriage returns keeps paragraph-mark bookkeeping consistent; the
# very first "`r" yields a placeholder with no run at all).
$r = $target.Range
$r.Collapse(0)
$r.InsertAfter("`r`r`r`r`r ")

$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# 1) Spacer paragraph
$p1 = $d.Paragraphs($startIdx + 1)
$p1.Range.InsertXML("<w:p $wns><w:pPr><w:spacing w:before='400'/></w:pPr></w:p>")

# 2) Red bar (top) - just two spaces, no explicit font
$p2 = $d.Paragraphs($startIdx + 2)
$p2.Range.InsertXML("<w:p $wns><w:pPr><w:shd w:fill='DC2626'/></w:pPr><w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>  </w:t></w:r></w:p>")
$p2b = $d.Paragraphs($startIdx + 2)
$p2b.Format.SpaceBefore = 0
$p2b.Format.SpaceAfter = 0

# 3) Heading: NOT GETTING INTERVIEW CALLBACKS?
$p3 = $d.Paragraphs($startIdx + 3)
$p3.Range.InsertXML("<w:p $wns><w:pPr><w:shd w:fill='DC2626'/><w:spacing w:before='200' w:after='100'/><w:jc w:val='center'/></w:pPr><w:r><w:rPr><w:rFonts w:ascii='Arial' w:cs='Arial' w:eastAsia='Arial' w:hAnsi='Arial'/><w:b/><w:bCs/><w:color w:val='FFFFFF'/><w:sz w:val='26'/><w:szCs w:val='26'/></w:rPr><w:t xml:space='preserve'>NOT GETTING INTERVIEW CALLBACKS?</w:t></w:r></w:p>")

# 4) Line: Most students make 3-5 critical CV mistakes without realising it.
$p4 = $d.Paragraphs($startIdx + 4)
$p4.Range.InsertXML("<w:p $wns><w:pPr><w:shd w:fill='DC2626'/><w:spacing w:after='100'/><w:jc w:val='center'/></w:pPr><w:r><w:rPr><w:rFonts w:ascii='Arial' w:cs='Arial' w:eastAsia='Arial' w:hAnsi='Arial'/><w:color w:val='FFFFFF'/><w:sz w:val='22'/><w:szCs w:val='22'/></w:rPr><w:t xml:space='preserve'>Most students make 3-5 critical CV mistakes without realising it.</w:t></w:r></w:p>")
$p4b = $d.Paragraphs($startIdx + 4)
$p4b.Format.SpaceBefore = 0

# 5) Line: Get a personalised video review of YOUR CV with specific fixes.
$p5 = $d.Paragraphs($startIdx + 5)
$p5.Range.InsertXML("<w:p $wns><w:pPr><w:shd w:fill='DC2626'/><w:spacing w:after='100'/><w:jc w:val='center'/></w:pPr><w:r><w:rPr><w:rFonts w:ascii='Arial' w:cs='Arial' w:eastAsia='Arial' w:hAnsi='Arial'/><w:color w:val='FFFFFF'/><w:sz w:val='22'/><w:szCs w:val='22'/></w:rPr><w:t xml:space='preserve'>Get a personalised video review of YOUR CV with specific fixes.</w:t></w:r></w:p>")
$p5b = $d.Paragraphs($startIdx + 5)
$p5b.Format.SpaceBefore = 0

# 6) CTA link: -> flyquest.co.za/cv
$p6 = $d.Paragraphs($startIdx + 6)
$p6.Range.InsertXML("<w:p $wns><w:pPr><w:shd w:fill='DC2626'/><w:spacing w:before='100' w:after='200'/><w:jc w:val='center'/></w:pPr><w:r><w:rPr><w:rFonts w:ascii='Arial' w:cs='Arial' w:eastAsia='Arial' w:hAnsi='Arial'/><w:b/><w:bCs/><w:color w:val='FFCF00'/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>" + [char]8594 + " flyquest.co.za/cv</w:t></w:r></w:p>")

Write-Host "Inserted CTA block after paragraph" $startIdx
